$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated site/cluster breakdown table (rows 2-9), columns:
#   A = site, B = cluster, C = n, D = clust_ratio (%)
$data = @(
    @("Cap Noir", 1, 5, 18.5185185185185),
    @("Cap Noir", 2, 9, 33.3333333333333),
    @("Cap Noir", 3, 9, 33.3333333333333),
    @("Cap Noir", 4, 4, 14.8148148148148),
    @("Pointe Suzanne", 1, 10, 32.258064516129),
    @("Pointe Suzanne", 2, 8, 25.8064516129032),
    @("Pointe Suzanne", 3, 6, 19.3548387096774),
    @("Pointe Suzanne", 4, 7, 22.5806451612903)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}
